$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The 10/18 (F18) and 10/22 (F22) plan entries are retyped as "復習" (review),
# and the React course plan entries that used to occupy F18:F21 (previously
# blocked by the now-removed Bootstrap course entries at F21:F22) are shifted
# up into F19:F21.
$ws.Range("F18").Value = "復習"
$ws.Range("F19").Value = "React系列课程从零基础到项目开发实战 44-55・練習"
$ws.Range("F20").Value = "React系列课程从零基础到项目开发实战 56-67・練習"
$ws.Range("F21").Value = "React系列课程从零基础到项目开发实战 68-83・練習"
$ws.Range("F22").Value = "復習"

# Mark 10/18 and 10/19 as completed (進捗 = 100%), matching the other
# finished days in the G column.
$ws.Range("G18").NumberFormat = "0%"
$ws.Range("G19").NumberFormat = "0%"
$ws.Range("G18").Value = 1
$ws.Range("G19").Value = 1

# Move the active selection to G20, where the user left off.
[void]$ws.Range("G20").Select()
